$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data (A1:D5) then write new table (A1:D8)
$ws.Range("A1:D8").ClearContents()

$headers = @("id", "multa", "valores", "dias")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @(1, 3, 23, 23),
    @(2, 4, 43, 1),
    @(3, 5, 54, 2),
    @(4, 6, 56, 23),
    @(5, 7, 54, 4),
    @(1, 8, 34, 32),
    @(2, 9, 32, 4)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws.Range("D8").Select()
